$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 663.2857
$ws.Range("I18").Value = 571.7273
$ws.Range("K18").Value = 571.7273
$ws.Range("M18").Value = -287.7273
$ws.Range("H34").Value = 4755.8335
$ws.Range("I34").Value = 2887
$ws.Range("K34").Value = 2887
$ws.Range("M34").Value = -2684
$ws.Range("H36").Value = 4755.8335
$ws.Range("I36").Value = 2887
$ws.Range("K36").Value = 2887
$ws.Range("M36").Value = -2172
$ws.Range("H62").Value = 2460.7144
$ws.Range("I62").Value = 2458.8
$ws.Range("K62").Value = 2458.8
$ws.Range("M62").Value = -1834.8
$ws.Range("H65").Value = 2460.7144
$ws.Range("I65").Value = 2458.8
$ws.Range("K65").Value = 12294
$ws.Range("M65").Value = -9174
$ws.Range("H112").Value = 21165410
$ws.Range("J112").Value = 22858562
$ws.Range("L112").Value = 68575686
$ws.Range("N112").Value = -68577902
$ws.Range("H139").Value = 49460
$ws.Range("J139").Value = 49460
$ws.Range("L139").Value = 49460
$ws.Range("N139").Value = -59740

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9758.219999999999
$ws.Range("I32").Value = 7215.7036
$ws.Range("J32").Value = 20597.37
$ws.Range("K32").Value = 7215.7036
$ws.Range("L32").Value = 20597.37
$ws.Range("M32").Value = -6928.7036
$ws.Range("N32").Value = -21171.37
$ws.Range("H74").Value = 8622001
$ws.Range("I74").Value = 1053.0555
$ws.Range("J74").Value = 22729008
$ws.Range("K74").Value = 1053.0555
$ws.Range("L74").Value = 22729008
$ws.Range("M74").Value = -179.0554999999999
$ws.Range("N74").Value = -22730756
$ws.Range("H77").Value = 8622001
$ws.Range("I77").Value = 1053.0555
$ws.Range("J77").Value = 22729008
$ws.Range("K77").Value = 5265.2775
$ws.Range("L77").Value = 113645040
$ws.Range("M77").Value = -897.2775000000001
$ws.Range("N77").Value = -113653776
$ws.Range("H132").Value = 1615608.2
$ws.Range("I132").Value = 1966.5264
$ws.Range("J132").Value = 4170541
$ws.Range("K132").Value = 5899.5792
$ws.Range("L132").Value = 12511623
$ws.Range("M132").Value = -3369.5792
$ws.Range("N132").Value = -12516683
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1005670
$ws.Range("I107").Value = 1507505.5
$ws.Range("J107").Value = 1999
$ws.Range("K107").Value = 1507505.5
$ws.Range("L107").Value = 1999
$ws.Range("M107").Value = -1505585.5
$ws.Range("N107").Value = -5839
$ws.Range("H134").Value = 38496.97
$ws.Range("I134").Value = 7103.92
$ws.Range("J134").Value = 169301.33
$ws.Range("K134").Value = 21311.76
$ws.Range("L134").Value = 507903.99
$ws.Range("M134").Value = -18776.76
$ws.Range("N134").Value = -512973.99

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21483772
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 21483772
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 21483772
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -21484362
$ws.Range("H34").Value = 21483772
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 21483772
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 21483772
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -21484176
$ws.Range("H132").Value = 3847985.5
$ws.Range("I132").Value = 5129349.5
$ws.Range("J132").Value = 3893.8462
$ws.Range("K132").Value = 15388048.5
$ws.Range("L132").Value = 11681.5386
$ws.Range("M132").Value = -15385518.5
$ws.Range("N132").Value = -16741.5386
$ws.Range("H135").Value = 38264.285
$ws.Range("J135").Value = 39284.617
$ws.Range("L135").Value = 39284.617
$ws.Range("N135").Value = -49424.617

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 3499.6667
$ws.Range("I76").Value = 2999
$ws.Range("J76").Value = 3750
$ws.Range("K76").Value = 8997
$ws.Range("L76").Value = 11250
$ws.Range("M76").Value = -8614
$ws.Range("N76").Value = -12016
$ws.Range("H79").Value = 3499.6667
$ws.Range("I79").Value = 2999
$ws.Range("J79").Value = 3750
$ws.Range("K79").Value = 8997
$ws.Range("L79").Value = 11250
$ws.Range("M79").Value = -7671
$ws.Range("N79").Value = -13902
$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 1500
$ws.Range("K82").Value = 4500
$ws.Range("M82").Value = -4094
$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 1500
$ws.Range("K85").Value = 4500
$ws.Range("M85").Value = -3096
$ws.Range("H100").Value = 3380
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""
$ws.Range("H109").Value = 1340.6364
$ws.Range("I109").Value = 963.8570999999999
$ws.Range("J109").Value = 2000
$ws.Range("K109").Value = 2891.5713
$ws.Range("L109").Value = 6000
$ws.Range("M109").Value = -1851.5713
$ws.Range("N109").Value = -8080
$ws.Range("H112").Value = 1998.0834
$ws.Range("I112").Value = 911
$ws.Range("J112").Value = 3520
$ws.Range("K112").Value = 2733
$ws.Range("L112").Value = 10560
$ws.Range("M112").Value = -1625
$ws.Range("N112").Value = -12776

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 45882.855
$ws.Range("J133").Value = 45882.855
$ws.Range("L133").Value = 45882.855
$ws.Range("N133").Value = -56002.855

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 45456572
$ws.Range("I68").Value = 1865.0834
$ws.Range("K68").Value = 1865.0834
$ws.Range("M68").Value = -1116.0834
$ws.Range("H71").Value = 45456572
$ws.Range("I71").Value = 1865.0834
$ws.Range("K71").Value = 9325.416999999999
$ws.Range("M71").Value = -5581.416999999999
$ws.Range("H136").Value = 9834.394
$ws.Range("I136").Value = 7945.9473
$ws.Range("J136").Value = 12397.286
$ws.Range("K136").Value = 23837.8419
$ws.Range("L136").Value = 37191.858
$ws.Range("M136").Value = -21287.8419
$ws.Range("N136").Value = -42291.858
$ws.Range("H137").Value = 54071.5
$ws.Range("J137").Value = 54885.8
$ws.Range("L137").Value = 54885.8
$ws.Range("N137").Value = -65085.8

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 131634.72
$ws.Range("I62").Value = 152573.83
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 152573.83
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -151949.83
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 131634.72
$ws.Range("I65").Value = 152573.83
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 762869.1499999999
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -759749.1499999999
$ws.Range("N65").Value = -36240
$ws.Range("H126").Value = 1620.1428
$ws.Range("I126").Value = 1468.2
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 4404.6
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -1934.6
$ws.Range("N126").Value = -10940
